$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '90.751.94'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +3.90%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.204.65'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.16%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '219.26'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +6.18%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '648.53'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +6.12%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.403'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +6.62%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.713'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +7.03%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '3.202.79'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.581'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +8.77%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000260'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +7.75%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.42'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.21%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '33.50'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +4.72%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '90.381.41'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.73%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.790.56'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +1.05%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.216.20'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +2.32%  '
$ws.Range('E19').Value = '  +11.68%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0000226'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +75.33%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '443.29'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +6.92%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '13.53'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.67'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.50%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.11'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.77%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.29'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.49%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.96'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '82.18'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +11.93%  '
$ws.Range('E28').Value = '  +1.25%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.161'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.40%  '
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.19'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +41.31%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '8.48'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.39%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '543.55'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '7.11'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +5.68%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.93'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +4.79%  '
$ws.Range('E37').Value = '  +0.68%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '22.59'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +3.78%  '
$ws.Range('E39').Value = '  +2.87%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.128'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.07%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('E42').Value = '  +2.73%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E44').Value = '  +1.73%  '
$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '147.21'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.96%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '45.22'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +4.77%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '173.71'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.761'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +9.21%  '
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.24'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.08%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.625'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +6.87%  '
